$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

$ws.Range("D2").Value = "29.308.53"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("E3").Value = "  -0.37%  "
Set-TextValue "D4" "0.9997"
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue "D5" "239.32"
$ws.Range("E5").Value = "  -0.40%  "
Set-TextValue "D6" "0.6249"
$ws.Range("E6").Value = "  -0.73%  "
Set-TextValue "D8" "0.07367"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.980.56"
$ws.Range("E11").Value = "  +7.27%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D12" "0.07717"
$ws.Range("E12").Value = "  -0.24%  "
Set-TextValue "D13" "4.944"
$ws.Range("E13").Value = "  -1.26%  "
Set-TextValue "D14" "0.00001061"
$ws.Range("E14").Value = "  +3.93%  "
Set-TextValue "D15" "0.6615"
$ws.Range("E15").Value = "  -2.92%  "
Set-TextValue "D16" "81.38"
$ws.Range("E16").Value = "  -1.23%  "
Set-TextValue "D17" "6.227"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "29.317.38"
$ws.Range("E18").Value = "  -0.63%  "
Set-TextValue "D19" "236.66"
$ws.Range("E19").Value = "  +3.19%  "
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("E22").Value = "  -3.65%  "
$ws.Range("E23").Value = "  +0.19%  "
Set-TextValue "D24" "157.27"
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("E25").Value = "  -1.11%  "
Set-TextValue "D26" "0.1333"
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("E27").Value = "  -1.43%  "
Set-TextValue "D28" "0.06985"
$ws.Range("E28").Value = "  +5.31%  "
Set-TextValue "D29" "1.471"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D31" "4.012"
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D32" "4.018"
$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("E33").Value = "  +1.01%  "
Set-TextValue "D34" "1.773"
$ws.Range("E34").Value = "  -3.91%  "
Set-TextValue "D35" "0.6782"
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("E36").Value = "  +0.41%  "
Set-TextValue "D37" "0.01822"
$ws.Range("E37").Value = "  -2.39%  "
Set-TextValue "D38" "2.776"
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("D39").Value = "1.230.13"
$ws.Range("E39").Value = "  -1.65%  "
Set-TextValue "D40" "6.715"
$ws.Range("E40").Value = "  -0.87%  "
Set-TextValue "D41" "0.9434"
$ws.Range("E41").Value = "  +0.49%  "
Set-TextValue "D42" "1.002"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "1.998.25"
$ws.Range("E43").Value = "  -0.77%  "
Set-TextValue "D44" "101.10"
$ws.Range("E44").Value = "  -0.09%  "
Set-TextValue "D45" "65.09"
$ws.Range("E45").Value = "  -1.51%  "
Set-TextValue "D46" "0.00000000120"
$ws.Range("E46").Value = "  +1.71%  "
Set-TextValue "D47" "6.923"
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("E48").Value = "  -2.09%  "
Set-TextValue "D49" "8.858"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("E50").Value = "  -2.49%  "
Set-TextValue "D51" "0.3869"
$ws.Range("E51").Value = "  -1.58%  "
